$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New product rows data: No, Name, Balance, OrderLimit, Price, SellPrice, Transactions
$rows = @(
    @{No=1; Name='AMARYL 1MG 30 TAB'; Balance='1:0'; Limit='1'; Price='40.00'; Sell='13.2000'; Trans='0:1'},
    @{No=2; Name='COLONA 30 F.C.TAB'; Balance='1:0'; Limit='1'; Price='69.00'; Sell='22.7700'; Trans='0:1'},
    @{No=3; Name='CONTROLOC 20MG 14  GASTRO RESISTANT TABS'; Balance='0:0'; Limit='1'; Price='188.00'; Sell='94.0000'; Trans='0:1'},
    @{No=4; Name='DAFLON 500MG 30 F.C. TABS'; Balance='0:0'; Limit='1'; Price='190.00'; Sell='95.0000'; Trans='0:1'},
    @{No=5; Name='DICYNONE 250MG/2ML 3 AMP.'; Balance='1:1'; Limit='1'; Price='42.00'; Sell='13.8600'; Trans='0:1'},
    @{No=6; Name='MEBEFAC 200 MG SR 30 F.C. TABS'; Balance='0:0'; Limit='1'; Price='66.00'; Sell='66.0000'; Trans='1:0'},
    @{No=7; Name='PRONTOGEST 100MG/2ML 10 IM AMPOULE'; Balance='1:3'; Limit='1'; Price='240.00'; Sell='24.0000'; Trans='0:1'},
    @{No=8; Name='سرنجات 3 سم'; Balance='0:0'; Limit='0'; Price='2.00'; Sell='4.0000'; Trans='2:0'}
)

# Insert 5 new rows before the current total row (old row 10) so existing rows 10 (total) and 11 (footer)
# shift down to 15 and 16, making room for product rows 10-14.
$ws.Rows.Item(10).Resize(5).Insert() | Out-Null

# Copy formatting for the new product rows (10-14) from the template row 9 (A9:Q9)
$ws.Range("A9:Q9").Copy() | Out-Null
$ws.Range("A10:Q14").PasteSpecial(-4122) | Out-Null
$ws.Range("A10:Q14").PasteSpecial(-4123) | Out-Null

# Set up merges for the new rows, matching the pattern used in row 7-9
for ($i = 0; $i -lt 5; $i++) {
    $r = 10 + $i
    $ws.Range("A$r`:B$r").Merge() | Out-Null
    $ws.Range("C$r`:G$r").Merge() | Out-Null
    $ws.Range("H$r`:K$r").Merge() | Out-Null
    $ws.Range("L$r`:M$r").Merge() | Out-Null
    $ws.Range("N$r`:O$r").Merge() | Out-Null
}

# Write the 8 product rows (rows 7 through 14)
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = 7 + $i
    $row = $rows[$i]
    $ws.Range("A$r").Value = $row.No
    $ws.Range("C$r").Value = $row.Name
    $ws.Range("H$r").Value = $row.Balance
    $ws.Range("L$r").Value = $row.Limit
    $ws.Range("N$r").Value = $row.Price
    $ws.Range("P$r").Value = $row.Sell
    $ws.Range("Q$r").Value = $row.Trans
}

# Row 15 (total) - formerly row 10
$ws.Range("P15").Value = 332.82999999999998

# Row 16 (footer) - formerly row 11
$ws.Range("A16").Value = "Saturday, 24 May, 2025 10:11 AM"
$ws.Range("G16").Value = "1/1"
$ws.Range("K16").Value = "developed by : Abdelaziz Talaat"
